$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status for the 2d28cac1...md row (row 3) changed from "Ready for handoff"
# to "Handback transform failed" everywhere that string is shown: the
# per-locale status sheets and the roll-up Overview sheet.
$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# New "Error Detail" (column K) values describing the handback mismatch.
$zhcn.Range("K3").Value = "Handback file name: nrgd1laz.vgj is different with handoff file name: 2d28cac1-9b54-4016-a77e-0dc5ee6706c6.0def3bfb436eeb5f67cb051b1aa1214b0b0658cb.zh-cn."
$dede.Range("K3").Value = "Handback file name: nrgd1laz.vgj is different with handoff file name: 2d28cac1-9b54-4016-a77e-0dc5ee6706c6.0def3bfb436eeb5f67cb051b1aa1214b0b0658cb.de-de."
